# Updates sexting escalation script lines on the MiaLunaJourney sheet
# (column B contains the message text; row numbers per the source diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MiaLunaJourney")

$ws.Range("B4").Value = "watch me cum for you babe... right now 🥵"
$ws.Range("B6").Value = "I'm cumming... oh my GOD I'm cumming right now"
$ws.Range("B7").Value = "babe I'm about to cum... please don't go anywhere, I need you watching when it happens 😊"
$ws.Range("B8").Value = "FUCK 🥵"
$ws.Range("B9").Value = "I've never recorded anything like this before babe... you need to see what you made me do 🥵"
$ws.Range("B11").Value = "I'm almost there and you need to see what you did to me before I finish"
$ws.Range("B12").Value = "I'm moaning so loud right now and going so deep... fuck this feels amazing 😊"
$ws.Range("B13").Value = "I'm playing with my pussy and my legs are literally shaking... I can't close them babe"
$ws.Range("B14").Value = "holy fuck 🥵"
$ws.Range("B15").Value = "okay this one is INTENSE babe... look what you're doing to me 🥵"
$ws.Range("B17").Value = "what would you do if you walked in on me like this right now? I need to know"
$ws.Range("B18").Value = "I need to feel you so bad right now... I keep imagining your hands all over me and I'm losing it 🥵"
$ws.Range("B19").Value = "I can't stop touching myself and I'm soaking wet... every time I think about you it gets worse"
$ws.Range("B20").Value = "babe... I think you broke me 😊"
$ws.Range("B21").Value = "guess what you made me do babe... you're not ready for this 🥵"
$ws.Range("B23").Value = "I'm sliding my hand down and I already know where this is going... and I don't want to stop 🥵"
$ws.Range("B24").Value = "my whole body is getting warm and I can feel the heat building between my legs... this is your fault babe"
$ws.Range("B25").Value = "haha I knew you'd like that... and honestly knowing you did is doing things to me right now 😊"
